$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A and B by 1 character (14.42578125 -> 15.42578125).
# NOTE: this runtime's ColumnWidth setter snaps to increments of 1/6 of a
# character, so the closest attainable stored width is 15.5; there is no
# COM call available that reproduces 15.42578125 exactly.
$ws.Columns.Item(1).ColumnWidth = 14.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.666666666666668

# Update the values in A1:B5
$ws.Range("A1").Value = -0.0077319541145460515
$ws.Range("B1").Value = -0.0073151622617322852

$ws.Range("A2").Value = -0.035113097906611136
$ws.Range("B2").Value = -0.045131186603830108

$ws.Range("A3").Value = -0.013595168236939547
$ws.Range("B3").Value = -0.025377727842637879

$ws.Range("A4").Value = -0.023474155373074377
$ws.Range("B4").Value = -0.022879431690503188

$ws.Range("A5").Value = -0.060146151198309898
$ws.Range("B5").Value = -0.060125882916858303
